$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("C2").Value = "тестовые данные"
$ws.Range("D2:F2").ClearContents()
$ws.Range("G2").Value = "[5..600]  "

# Row 3 updates
$ws.Range("B3").Value = "Не найден  "
$ws.Range("C3:F3").ClearContents()
